$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 18
$ws_ALC.Range("H18").Value = 19233844
$ws_ALC.Range("I18").Value = 19233844
$ws_ALC.Range("K18").Value = 19233844
$ws_ALC.Range("M18").Value = -19233560

# ALC row 116
$ws_ALC.Range("H116").Value = 6878.5884
$ws_ALC.Range("I116").Value = 4692.143
$ws_ALC.Range("K116").Value = 4692.143
$ws_ALC.Range("M116").Value = -1250.143

# ALC row 138
$ws_ALC.Range("H138").Value = 2890.2444
$ws_ALC.Range("I138").Value = 1704.1111
$ws_ALC.Range("J138").Value = 3681
$ws_ALC.Range("K138").Value = 5112.3333
$ws_ALC.Range("L138").Value = 11043
$ws_ALC.Range("M138").Value = 27.66669999999976
$ws_ALC.Range("N138").Value = -21323

# ARM row 32
$ws_ARM.Range("H32").Value = 8639.166999999999
$ws_ARM.Range("I32").Value = 4881.4814
$ws_ARM.Range("J32").Value = 19912.223
$ws_ARM.Range("K32").Value = 4881.4814
$ws_ARM.Range("L32").Value = 19912.223
$ws_ARM.Range("M32").Value = -4594.4814
$ws_ARM.Range("N32").Value = -20486.223

# ARM row 43
$ws_ARM.Range("H43").Value = 16292.167
$ws_ARM.Range("J43").Value = 20938.5
$ws_ARM.Range("L43").Value = 20938.5
$ws_ARM.Range("N43").Value = -21564.5

# ARM row 45
$ws_ARM.Range("H45").Value = 6543173.5
$ws_ARM.Range("I45").Value = 15986006
$ws_ARM.Range("K45").Value = 15986006
$ws_ARM.Range("M45").Value = -15985629

# ARM row 61
$ws_ARM.Range("H61").Value = 6674.8
$ws_ARM.Range("I61").Value = 6861.5
$ws_ARM.Range("J61").Value = 4994.5
$ws_ARM.Range("K61").Value = 6861.5
$ws_ARM.Range("L61").Value = 4994.5
$ws_ARM.Range("M61").Value = -6649.5
$ws_ARM.Range("N61").Value = -5418.5

# ARM row 63
$ws_ARM.Range("H63").Value = 4475.304
$ws_ARM.Range("I63").Value = 1654.091
$ws_ARM.Range("J63").Value = 7061.4165
$ws_ARM.Range("K63").Value = 1654.091
$ws_ARM.Range("L63").Value = 7061.4165
$ws_ARM.Range("M63").Value = -968.0909999999999
$ws_ARM.Range("N63").Value = -8433.416499999999

# ARM row 66
$ws_ARM.Range("H66").Value = 4475.304
$ws_ARM.Range("I66").Value = 1654.091
$ws_ARM.Range("J66").Value = 7061.4165
$ws_ARM.Range("K66").Value = 8270.455
$ws_ARM.Range("L66").Value = 35307.0825
$ws_ARM.Range("M66").Value = -4838.455
$ws_ARM.Range("N66").Value = -42171.0825

# ARM row 74
$ws_ARM.Range("H74").Value = 26042.625
$ws_ARM.Range("I74").Value = 3633.9722
$ws_ARM.Range("K74").Value = 3633.9722
$ws_ARM.Range("M74").Value = -2759.9722

# ARM row 77
$ws_ARM.Range("H77").Value = 26042.625
$ws_ARM.Range("I77").Value = 3633.9722
$ws_ARM.Range("K77").Value = 18169.861
$ws_ARM.Range("M77").Value = -13801.861

# ARM row 110
$ws_ARM.Range("H110").Value = 927010.4399999999
$ws_ARM.Range("I110").Value = 1029863.5
$ws_ARM.Range("J110").Value = 1333
$ws_ARM.Range("K110").Value = 1029863.5
$ws_ARM.Range("L110").Value = 1333
$ws_ARM.Range("M110").Value = -1027818.5
$ws_ARM.Range("N110").Value = -5423

# ARM row 136
$ws_ARM.Range("H136").Value = 6674.8
$ws_ARM.Range("I136").Value = 6861.5
$ws_ARM.Range("J136").Value = 4994.5
$ws_ARM.Range("K136").Value = 20584.5
$ws_ARM.Range("L136").Value = 14983.5
$ws_ARM.Range("M136").Value = -18034.5
$ws_ARM.Range("N136").Value = -20083.5

# BSM row 20
$ws_BSM.Range("H20").Value = 37044184
$ws_BSM.Range("I20").Value = 41673460
$ws_BSM.Range("K20").Value = 41673460
$ws_BSM.Range("M20").Value = -41673213

# BSM row 42
$ws_BSM.Range("H42").Value = 285001
$ws_BSM.Range("J42").Value = 285001
$ws_BSM.Range("L42").Value = 285001
$ws_BSM.Range("N42").Value = -285657

# BSM row 99
$ws_BSM.Range("H99").Value = 17859578
$ws_BSM.Range("I99").Value = 20410374
$ws_BSM.Range("K99").Value = 20410374
$ws_BSM.Range("M99").Value = -20408876

# BSM row 105
$ws_BSM.Range("H105").Value = 5682831.5
$ws_BSM.Range("I105").Value = 5682831.5
$ws_BSM.Range("J105").Value = 0
$ws_BSM.Range("K105").Value = 5682831.5
$ws_BSM.Range("L105").Value = 0
$ws_BSM.Range("M105").Value = -5681084.5
$ws_BSM.Range("N105").ClearContents()

# BSM row 109
$ws_BSM.Range("H109").Value = 79995
$ws_BSM.Range("J109").Value = 79995
$ws_BSM.Range("L109").Value = 79995
$ws_BSM.Range("N109").Value = -82769

# BSM row 134
$ws_BSM.Range("H134").Value = 3786.4666
$ws_BSM.Range("I134").Value = 1299.7727
$ws_BSM.Range("J134").Value = 10624.875
$ws_BSM.Range("K134").Value = 3899.3181
$ws_BSM.Range("L134").Value = 31874.625
$ws_BSM.Range("M134").Value = -1364.3181
$ws_BSM.Range("N134").Value = -36944.625

# CRP row 31
$ws_CRP.Range("H31").Value = 23291.404
$ws_CRP.Range("I31").Value = 2982.3809
$ws_CRP.Range("J31").Value = 39694.848
$ws_CRP.Range("K31").Value = 2982.3809
$ws_CRP.Range("L31").Value = 39694.848
$ws_CRP.Range("M31").Value = -2687.3809
$ws_CRP.Range("N31").Value = -40284.848

# CRP row 34
$ws_CRP.Range("H34").Value = 23291.404
$ws_CRP.Range("I34").Value = 2982.3809
$ws_CRP.Range("J34").Value = 39694.848
$ws_CRP.Range("K34").Value = 2982.3809
$ws_CRP.Range("L34").Value = 39694.848
$ws_CRP.Range("M34").Value = -2780.3809
$ws_CRP.Range("N34").Value = -40098.848

# CRP row 58
$ws_CRP.Range("H58").Value = 7302.4346
$ws_CRP.Range("I58").Value = 11564.8
$ws_CRP.Range("K58").Value = 11564.8
$ws_CRP.Range("M58").Value = -11361.8

# CRP row 99
$ws_CRP.Range("H99").Value = 3611.5652
$ws_CRP.Range("I99").Value = 3550.0715
$ws_CRP.Range("J99").Value = 3707.2222
$ws_CRP.Range("K99").Value = 3550.0715
$ws_CRP.Range("L99").Value = 3707.2222
$ws_CRP.Range("M99").Value = -2052.0715
$ws_CRP.Range("N99").Value = -6703.2222

# CRP row 122
$ws_CRP.Range("H122").Value = 3434.75
$ws_CRP.Range("I122").Value = 3664.7144
$ws_CRP.Range("J122").Value = 1825
$ws_CRP.Range("K122").Value = 10994.1432
$ws_CRP.Range("L122").Value = 5475
$ws_CRP.Range("M122").Value = -8544.143199999999
$ws_CRP.Range("N122").Value = -10375

# CRP row 125
$ws_CRP.Range("H125").Value = 0
$ws_CRP.Range("J125").Value = 0
$ws_CRP.Range("L125").Value = 0
$ws_CRP.Range("N125").ClearContents()

# CRP row 126
$ws_CRP.Range("H126").Value = 3611.5652
$ws_CRP.Range("I126").Value = 3550.0715
$ws_CRP.Range("J126").Value = 3707.2222
$ws_CRP.Range("K126").Value = 10650.2145
$ws_CRP.Range("L126").Value = 11121.6666
$ws_CRP.Range("M126").Value = -8180.2145
$ws_CRP.Range("N126").Value = -16061.6666

# CRP row 127
$ws_CRP.Range("H127").Value = 54898.332
$ws_CRP.Range("J127").Value = 54898.332
$ws_CRP.Range("L127").Value = 54898.332
$ws_CRP.Range("N127").Value = -64818.332

# CRP row 131
$ws_CRP.Range("H131").Value = 99500
$ws_CRP.Range("J131").Value = 99500
$ws_CRP.Range("L131").Value = 99500
$ws_CRP.Range("N131").Value = -109580

# CRP row 132
$ws_CRP.Range("H132").Value = 46940.69
$ws_CRP.Range("I132").Value = 33782.13
$ws_CRP.Range("J132").Value = 84023.91
$ws_CRP.Range("K132").Value = 101346.39
$ws_CRP.Range("L132").Value = 252071.73
$ws_CRP.Range("M132").Value = -98816.38999999998
$ws_CRP.Range("N132").Value = -257131.73

# CRP row 134
$ws_CRP.Range("H134").Value = 2471.756
$ws_CRP.Range("I134").Value = 1473.742
$ws_CRP.Range("K134").Value = 4421.226
$ws_CRP.Range("M134").Value = -1886.226

# CRP row 135
$ws_CRP.Range("H135").Value = 99060.55499999999
$ws_CRP.Range("J135").Value = 99060.55499999999
$ws_CRP.Range("L135").Value = 99060.55499999999
$ws_CRP.Range("N135").Value = -109200.555

# CRP row 136
$ws_CRP.Range("H136").Value = 7302.4346
$ws_CRP.Range("I136").Value = 11564.8
$ws_CRP.Range("K136").Value = 34694.39999999999
$ws_CRP.Range("M136").Value = -32144.39999999999

# CRP row 138
$ws_CRP.Range("H138").Value = 89998.5
$ws_CRP.Range("J138").Value = 89998.5
$ws_CRP.Range("L138").Value = 89998.5
$ws_CRP.Range("N138").Value = -100278.5

# CRP row 140
$ws_CRP.Range("H140").Value = 120000
$ws_CRP.Range("J140").Value = 120000
$ws_CRP.Range("L140").Value = 120000
$ws_CRP.Range("N140").Value = -130360

# CRP row 141
$ws_CRP.Range("H141").Value = 196911.17
$ws_CRP.Range("J141").Value = 212994
$ws_CRP.Range("L141").Value = 212994
$ws_CRP.Range("N141").Value = -223354

# CUL row 44
$ws_CUL.Range("H44").Value = 4042.8572
$ws_CUL.Range("I44").Value = 1766.6666
$ws_CUL.Range("K44").Value = 5299.9998
$ws_CUL.Range("M44").Value = -4901.9998

# CUL row 103
$ws_CUL.Range("H103").Value = 370.33334
$ws_CUL.Range("I103").Value = 344.4
$ws_CUL.Range("K103").Value = 1033.2
$ws_CUL.Range("M103").Value = -154.1999999999998

# GSM row 70
$ws_GSM.Range("H70").Value = 50004124
$ws_GSM.Range("I70").Value = 100002750
$ws_GSM.Range("K70").Value = 100002750
$ws_GSM.Range("M70").Value = -100002480

# GSM row 73
$ws_GSM.Range("H73").Value = 50004124
$ws_GSM.Range("I73").Value = 100002750
$ws_GSM.Range("K73").Value = 100002750
$ws_GSM.Range("M73").Value = -100001814

# GSM row 97
$ws_GSM.Range("H97").Value = 1701099.9
$ws_GSM.Range("J97").Value = 245
$ws_GSM.Range("L97").Value = 245
$ws_GSM.Range("N97").Value = -1237

# GSM row 102
$ws_GSM.Range("H102").Value = 6427444
$ws_GSM.Range("I102").Value = 10102164
$ws_GSM.Range("K102").Value = 10102164
$ws_GSM.Range("M102").Value = -10100542

# GSM row 122
$ws_GSM.Range("H122").Value = 310151.4
$ws_GSM.Range("I122").Value = 406801.5
$ws_GSM.Range("J122").Value = 6394
$ws_GSM.Range("K122").Value = 1220404.5
$ws_GSM.Range("L122").Value = 19182
$ws_GSM.Range("M122").Value = -1217954.5
$ws_GSM.Range("N122").Value = -24082

# GSM row 132
$ws_GSM.Range("H132").Value = 3602.68
$ws_GSM.Range("I132").Value = 3544.6667
$ws_GSM.Range("J132").Value = 4995
$ws_GSM.Range("K132").Value = 10634.0001
$ws_GSM.Range("L132").Value = 14985
$ws_GSM.Range("M132").Value = -8104.000100000001
$ws_GSM.Range("N132").Value = -20045

# LTW row 7
$ws_LTW.Range("H7").Value = 11453
$ws_LTW.Range("I7").Value = 9866.5
$ws_LTW.Range("K7").Value = 9866.5
$ws_LTW.Range("M7").Value = -9754.5

# LTW row 22
$ws_LTW.Range("H22").Value = 20825.305
$ws_LTW.Range("I22").Value = 60164.934
$ws_LTW.Range("J22").Value = 1790
$ws_LTW.Range("K22").Value = 60164.934
$ws_LTW.Range("L22").Value = 1790
$ws_LTW.Range("M22").Value = -59869.934
$ws_LTW.Range("N22").Value = -2380

# LTW row 27
$ws_LTW.Range("H27").Value = 20825.305
$ws_LTW.Range("I27").Value = 60164.934
$ws_LTW.Range("J27").Value = 1790
$ws_LTW.Range("K27").Value = 60164.934
$ws_LTW.Range("L27").Value = 1790
$ws_LTW.Range("M27").Value = -60057.934
$ws_LTW.Range("N27").Value = -2004

# LTW row 39
$ws_LTW.Range("H39").Value = 14000
$ws_LTW.Range("J39").Value = 14000
$ws_LTW.Range("L39").Value = 14000
$ws_LTW.Range("N39").Value = -14920

# LTW row 46
$ws_LTW.Range("H46").Value = 1179502.4
$ws_LTW.Range("I46").Value = 10870010
$ws_LTW.Range("K46").Value = 10870010
$ws_LTW.Range("M46").Value = -10869822

# LTW row 100
$ws_LTW.Range("H100").Value = 2361.1614
$ws_LTW.Range("J100").Value = 2584.1667
$ws_LTW.Range("L100").Value = 2584.1667
$ws_LTW.Range("N100").Value = -3666.1667

# LTW row 104
$ws_LTW.Range("H104").Value = 27188.2
$ws_LTW.Range("J104").Value = 27188.2
$ws_LTW.Range("L104").Value = 27188.2
$ws_LTW.Range("N104").Value = -34176.2

# LTW row 126
$ws_LTW.Range("H126").Value = 11453
$ws_LTW.Range("I126").Value = 9866.5
$ws_LTW.Range("K126").Value = 29599.5
$ws_LTW.Range("M126").Value = -27129.5

# LTW row 132
$ws_LTW.Range("H132").Value = 7799.4736
$ws_LTW.Range("I132").Value = 7920.636
$ws_LTW.Range("K132").Value = 23761.908
$ws_LTW.Range("M132").Value = -21231.908

# WVR row 32
$ws_WVR.Range("H32").Value = 9995
$ws_WVR.Range("I32").Value = 9995
$ws_WVR.Range("J32").Value = 0
$ws_WVR.Range("K32").Value = 9995
$ws_WVR.Range("L32").Value = 0
$ws_WVR.Range("N32").Value = -9678
$ws_WVR.Range("M32").ClearContents()

# WVR row 75
$ws_WVR.Range("H75").Value = 55000
$ws_WVR.Range("I75").Value = 55000
$ws_WVR.Range("K75").Value = 55000
$ws_WVR.Range("M75").Value = -54064

# WVR row 78
$ws_WVR.Range("H78").Value = 55000
$ws_WVR.Range("I78").Value = 55000
$ws_WVR.Range("K78").Value = 165000
$ws_WVR.Range("M78").Value = -160320

# WVR row 113
$ws_WVR.Range("H113").Value = 1516.4572
$ws_WVR.Range("J113").Value = 2299.875
$ws_WVR.Range("L113").Value = 6899.625
$ws_WVR.Range("N113").Value = -11239.625

# WVR row 123
$ws_WVR.Range("H123").Value = 98279
$ws_WVR.Range("J123").Value = 98279
$ws_WVR.Range("L123").Value = 98279
$ws_WVR.Range("N123").Value = -108079

# WVR row 132
$ws_WVR.Range("H132").Value = 21064316
$ws_WVR.Range("I132").Value = 27782642
$ws_WVR.Range("K132").Value = 83347926
$ws_WVR.Range("M132").Value = -83345396

# WVR row 136
$ws_WVR.Range("H136").Value = 6650.1665
$ws_WVR.Range("I136").Value = 6940.2
$ws_WVR.Range("K136").Value = 20820.6
$ws_WVR.Range("M136").Value = -18270.6
